$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.2651245551601423
$ws1.Range("C2").Value = 0.06349206349206349
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 0.1194029850746269
$ws1.Range("F2").Value = 0.2531645569620253
$ws1.Range("G2").Value = 0.6380368098159509
$ws1.Range("H2").Value = 0.7780230069555912
$ws1.Range("I2").Value = 28
$ws1.Range("J2").Value = 413
$ws1.Range("K2").Value = 121
$ws1.Range("L2").Value = 0

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0.2265917602996255
$ws2.Range("D2").Value = 0.3694656488549619

$ws2.Range("B3").Value = 0.06349206349206349
$ws2.Range("C3").Value = 1
$ws2.Range("D3").Value = 0.1194029850746269

$ws2.Range("B4").Value = 0.2651245551601423
$ws2.Range("C4").Value = 0.2651245551601423
$ws2.Range("D4").Value = 0.2651245551601423
$ws2.Range("E4").Value = 0.2651245551601423

$ws2.Range("B5").Value = 0.5317460317460317
$ws2.Range("C5").Value = 0.6132958801498127
$ws2.Range("D5").Value = 0.2444343169647944

$ws2.Range("B6").Value = 0.9533412415974695
$ws2.Range("C6").Value = 0.2651245551601423
$ws2.Range("D6").Value = 0.3570070108018491

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 121
$ws3.Range("C2").Value = 413
$ws3.Range("B3").Value = 0
$ws3.Range("C3").Value = 28
